# Update PERMANOVA p-values in column F to rounded (3-digit) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 0.5874
    "F3"  = 0.0018
    "F4"  = 0.6864
    "F7"  = 0.348
    "F8"  = 0.8598
    "F9"  = 0.6089
    "F12" = 0.5362
    "F13" = 0.947
    "F14" = 0.2406
    "F17" = 0.7767
    "F18" = 0.3417
    "F19" = 0.8712
    "F22" = 0.1537
    "F23" = 0.096
    "F24" = 0.2412
    "F27" = 0.2568
    "F28" = 0.6478
    "F29" = 0.7012
    "F32" = 0.876
    "F33" = 0.0004
    "F34" = 0.9858
    "F37" = 0.3882
    "F38" = 0.043
    "F39" = 0.4517
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
